$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("stimFile") before the existing cs_plus/cs_minus columns,
# which shifts those two columns to C/D.
$ws.Range("B1").EntireColumn.Insert()

# New header row (row 1)
$ws.Range("A1").Value = "posFile"
$ws.Range("B1").Value = "stimFile"
$ws.Range("C1").Value = "cs_plus_s"
$ws.Range("D1").Value = "cs_minus_s"
$ws.Range("E1").Value = "cs_plus_ns"
$ws.Range("F1").Value = "cs_minus_ns"

# New data row (row 2) replaces the old social row (row 2) and merges in the
# non-social columns that used to live on row 3.
$ws.Range("A2").Value = "positions.xlsx"
$ws.Range("B2").Value = "stimuli.xlsx"
$ws.Range("C2").Value = "stimuli/social/031_y_m_n_a.jpg"
$ws.Range("D2").Value = "stimuli/social/016_y_m_n_b.jpg"
$ws.Range("E2").Value = "stimuli/non-social/031_y_m_n_a_scrambled.jpg"
$ws.Range("F2").Value = "stimuli/non-social/016_y_m_n_b_scrambled.jpg"

# Old row 3 (previously the non-social row) is no longer needed now that it is
# merged into row 2.
$ws.Range("A3:F3").Delete()

# Column B should be the same visual width as column A (posFile), but not
# marked as auto bestFit (it was manually sized instead).
$ws.Cells.Item(1, 2).ColumnWidth = 15.63

# Match the selection shown after the edit.
$ws.Range("B1:B2").Select()
